# Updated with WorkFlow Execution
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns: StoreResponseVariables (L) / AddifyVariables (M) ---
$ws.Range("L1").Value = "StoreResponseVariables"
$ws.Range("M1").Value = "AddifyVariables"

# --- Row 2 becomes the "EDI-271 API test" entry (TestCaseName cleared) ---
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "EDI-271 API test"
$ws.Range("C2").Value = "http://mockbin.org/bin/3f64e65d-c657-42d5-bcc9-5b13e71ca493"

# --- Row 3 (PetGet) gains a StoreResponseVariables value ---
$ws.Range("L3").Value = "petId=id;petName=name"

# --- New row 4: PetPost test case now hitting the live API ---
$ws.Range("A4").Value = "PetPost"
$ws.Range("B4").Value = "post API Testing"
$ws.Range("C4").Value = "https://live.virtualandemo.com/api/pets"
$ws.Range("D4").Value = "application/json"
$ws.Range("E4").Value = "post-request.json"
$ws.Range("G4").Value = "post-response.json"
$ws.Range("I4").Value = "POST"
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = "Id=[petId];petName=doggie"

# --- Hyperlinks: update existing C2 link target, add new one for C4 ---
$ws.Hyperlinks.Item(1).Address = "http://mockbin.org/bin/3f64e65d-c657-42d5-bcc9-5b13e71ca493"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://live.virtualandemo.com/api/pets")
